$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '28.218.02'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  +3.46%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.915.22'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  +2.80%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.006'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  -1.26%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '316.18'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +1.08%  '
$ws.Range('E6').Value = '  -1.19%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.4860'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +1.22%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3841'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +2.94%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.07419'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -0.05%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.9556'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +1.83%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '21.03'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +1.28%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.07826'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -0.64%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '1.901.80'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +1.80%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '5.572'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +2.49%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '6.668'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +1.89%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '92.45'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +2.34%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '1.008'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -1.19%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.000008932'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +1.79%  '
$ws.Range('E19').Value = '  -1.17%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '28.236.23'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +3.37%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '15.10'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +2.01%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '5.188'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +1.30%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '2.146.10'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +1.99%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '10.95'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +2.40%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '1.963'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +0.53%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '157.03'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +1.84%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '18.65'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +0.59%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '2.122'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +5.69%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '117.09'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +0.94%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '5.047'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +0.89%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.08924'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +0.12%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '3.340'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -0.15%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '1.255'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +4.60%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.7852'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +5.62%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '4.718'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +3.20%  '
$ws.Range('E36').Value = '  +3.45%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '1.136'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +0.96%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.02057'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +0.35%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.05407'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +2.20%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.5611'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +4.49%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '3.027'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +1.06%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '7.152'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +0.46%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '8.643'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +2.84%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.1540'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +0.04%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.4973'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +3.38%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '10.82'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +1.86%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '107.83'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +4.53%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '1.686'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +1.81%  '
$ws.Range('E49').Value = '  -1.27%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '69.70'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +4.41%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.06136'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +0.89%  '
